# Add the "Crystal Lakes and Sourdough Gap" hike as a new row to Table1
# on the "Hike Difficulties" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing structured table and append a new row to it so the
# table range / autofilter / dimension all expand together (matches the
# table ref going from A1:D29 -> A1:D30 in the saved file).
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$rowRange = $newRow.Range
$rowRange.Item(1, 1).Value = "Crystal Lakes and Sourdough Gap"
$rowRange.Item(1, 2).Value = 7.5
$rowRange.Item(1, 3).Value = 2970
$rowRange.Item(1, 4).Value = "moderate"

# Match the post-edit selection (the cell just below the newly added row).
$ws.Range("D31").Select()
